$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 48616236
$ws.Range("I62").Value = 20005530
$ws.Range("K62").Value = 20005530
$ws.Range("M62").Value = -20004906

$ws.Range("H65").Value = 48616236
$ws.Range("I65").Value = 20005530
$ws.Range("K65").Value = 100027650
$ws.Range("M65").Value = -100024530

$ws.Range("H74").Value = 3972.2778
$ws.Range("I74").Value = 4714.2856
$ws.Range("J74").Value = 3500.0908
$ws.Range("K74").Value = 4714.2856
$ws.Range("L74").Value = 3500.0908
$ws.Range("M74").Value = -3778.2856
$ws.Range("N74").Value = -5372.0908

$ws.Range("H77").Value = 3972.2778
$ws.Range("I77").Value = 4714.2856
$ws.Range("J77").Value = 3500.0908
$ws.Range("K77").Value = 23571.428
$ws.Range("L77").Value = 17500.454
$ws.Range("M77").Value = -18891.428
$ws.Range("N77").Value = -26860.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("K21").Value = 1000
$ws.Range("M21").Value = -626

$ws.Range("H32").Value = 6726104
$ws.Range("I32").Value = 1627875.1
$ws.Range("K32").Value = 1627875.1
$ws.Range("M32").Value = -1627588.1

$ws.Range("H45").Value = 345848.25
$ws.Range("I45").Value = 667500.7
$ws.Range("J45").Value = 1220.6428
$ws.Range("K45").Value = 667500.7
$ws.Range("L45").Value = 1220.6428
$ws.Range("M45").Value = -667123.7
$ws.Range("N45").Value = -1974.6428

$ws.Range("H61").Value = 3866178
$ws.Range("I61").Value = 1985385
$ws.Range("K61").Value = 1985385
$ws.Range("M61").Value = -1985173

$ws.Range("H97").Value = 751.2857
$ws.Range("I97").Value = 1500
$ws.Range("J97").Value = 626.5
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 626.5
$ws.Range("M97").Value = -1004
$ws.Range("N97").Value = -1618.5

$ws.Range("H102").Value = 2670
$ws.Range("I102").Value = 2812.5
$ws.Range("J102").Value = 2100
$ws.Range("K102").Value = 2812.5
$ws.Range("L102").Value = 2100
$ws.Range("M102").Value = -1190.5
$ws.Range("N102").Value = -5344

$ws.Range("H122").Value = 1353.3684
$ws.Range("I122").Value = 697.5
$ws.Range("J122").Value = 1830.3636
$ws.Range("K122").Value = 2092.5
$ws.Range("L122").Value = 5491.0908
$ws.Range("M122").Value = 357.5
$ws.Range("N122").Value = -10391.0908

$ws.Range("H136").Value = 3866178
$ws.Range("I136").Value = 1985385
$ws.Range("K136").Value = 5956155
$ws.Range("M136").Value = -5953605

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1828.7407
$ws.Range("I94").Value = 1441.125
$ws.Range("J94").Value = 2392.5454
$ws.Range("K94").Value = 1441.125
$ws.Range("L94").Value = 2392.5454
$ws.Range("M94").Value = -990.125
$ws.Range("N94").Value = -3294.5454

$ws.Range("H105").Value = 1665.5555
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 1398
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 1398
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -4892

$ws.Range("H134").Value = 24352668
$ws.Range("I134").Value = 38463060
$ws.Range("J134").Value = 3970988.8
$ws.Range("K134").Value = 115389180
$ws.Range("L134").Value = 11912966.4
$ws.Range("M134").Value = -115386645
$ws.Range("N134").Value = -11918036.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3532650.5
$ws.Range("I58").Value = 1988245.5
$ws.Range("K58").Value = 1988245.5
$ws.Range("M58").Value = -1988042.5

$ws.Range("H134").Value = 2677000.2
$ws.Range("I134").Value = 12882
$ws.Range("K134").Value = 38646
$ws.Range("M134").Value = -36111

$ws.Range("H136").Value = 3532650.5
$ws.Range("I136").Value = 1988245.5
$ws.Range("K136").Value = 5964736.5
$ws.Range("M136").Value = -5962186.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 76
$ws.Range("I26").Value = 76
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 228
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 60
$ws.Range("N26").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5376
$ws.Range("I102").Value = 5867.4
$ws.Range("K102").Value = 5867.4
$ws.Range("M102").Value = -4245.4

$ws.Range("H132").Value = 14720496
$ws.Range("I132").Value = 20635802
$ws.Range("K132").Value = 61907406
$ws.Range("M132").Value = -61904876

$ws.Range("H134").Value = 24163
$ws.Range("J134").Value = 24163
$ws.Range("L134").Value = 72489
$ws.Range("N134").Value = -77559

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1261.0769
$ws.Range("I7").Value = 867.2
$ws.Range("J7").Value = 1798.1818
$ws.Range("K7").Value = 867.2
$ws.Range("L7").Value = 1798.1818
$ws.Range("M7").Value = -755.2
$ws.Range("N7").Value = -2022.1818

$ws.Range("H40").Value = 2640.3333
$ws.Range("I40").Value = 1725
$ws.Range("J40").Value = 3098
$ws.Range("K40").Value = 1725
$ws.Range("L40").Value = 3098
$ws.Range("M40").Value = -1589
$ws.Range("N40").Value = -3370

$ws.Range("H46").Value = 378.27274
$ws.Range("J46").Value = 352.33334
$ws.Range("L46").Value = 352.33334
$ws.Range("N46").Value = -728.33334

$ws.Range("H68").Value = 2477.2666
$ws.Range("J68").Value = 2475.6428
$ws.Range("L68").Value = 2475.6428
$ws.Range("N68").Value = -3973.6428

$ws.Range("H71").Value = 2477.2666
$ws.Range("J71").Value = 2475.6428
$ws.Range("L71").Value = 12378.214
$ws.Range("N71").Value = -19866.214

$ws.Range("H93").Value = 11898.458
$ws.Range("I93").Value = 3604.5
$ws.Range("J93").Value = 17822.715
$ws.Range("K93").Value = 3604.5
$ws.Range("L93").Value = 17822.715
$ws.Range("M93").Value = -2356.5
$ws.Range("N93").Value = -20318.715

$ws.Range("H122").Value = 9623217
$ws.Range("I122").Value = 1332317.5
$ws.Range("J122").Value = 28573844
$ws.Range("K122").Value = 3996952.5
$ws.Range("L122").Value = 85721532
$ws.Range("M122").Value = -3994502.5
$ws.Range("N122").Value = -85726432

$ws.Range("H126").Value = 1261.0769
$ws.Range("I126").Value = 867.2
$ws.Range("J126").Value = 1798.1818
$ws.Range("K126").Value = 2601.6
$ws.Range("L126").Value = 5394.5454
$ws.Range("M126").Value = -131.6000000000004
$ws.Range("N126").Value = -10334.5454

$ws.Range("H136").Value = 11766449
$ws.Range("I136").Value = 23530958
$ws.Range("J136").Value = 1940.5454
$ws.Range("K136").Value = 70592874
$ws.Range("L136").Value = 5821.6362
$ws.Range("M136").Value = -70590324
$ws.Range("N136").Value = -10921.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1343.2222
$ws.Range("I96").Value = 1469.5
$ws.Range("J96").Value = 1090.6666
$ws.Range("K96").Value = 1469.5
$ws.Range("L96").Value = 1090.6666
$ws.Range("M96").Value = -96.5
$ws.Range("N96").Value = -3836.6666

$ws.Range("H100").Value = 7853.077
$ws.Range("I100").Value = 8482.5
$ws.Range("K100").Value = 16965
$ws.Range("M100").Value = -16424

$ws.Range("H122").Value = 1147.6666
$ws.Range("I122").Value = 1017.2
$ws.Range("J122").Value = 1444.1818
$ws.Range("K122").Value = 3051.6
$ws.Range("L122").Value = 4332.5454
$ws.Range("M122").Value = -601.6000000000004
$ws.Range("N122").Value = -9232.545399999999

$ws.Range("H126").Value = 20835246
$ws.Range("I126").Value = 25000656
$ws.Range("J126").Value = 8199.5
$ws.Range("K126").Value = 75001968
$ws.Range("L126").Value = 24598.5
$ws.Range("M126").Value = -74999498
$ws.Range("N126").Value = -29538.5
